$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MS")

# Update the account labels in column A with their new trailing account numbers.
$ws.Range("A3").Value  = "Cardor Alef SpA - 4597"
$ws.Range("A4").Value  = "Cuenta Personal - Jaime - 4231"
$ws.Range("A5").Value  = "NNW Capital SpA - 4757"
$ws.Range("A6").Value  = "NNW II Capital SpA - 4751"
$ws.Range("A7").Value  = "FNW Capital SpA - 4756"
$ws.Range("A10").Value = "Alanseb LP - 4582"
$ws.Range("A11").Value = "Cuenta Personal - Felipe - 4228"
$ws.Range("A12").Value = "Cuenta Personal - Irene - 4441"
$ws.Range("A13").Value = "Cuenta Personal - Jaime y Felipe - 4191"
$ws.Range("A14").Value = "Cuenta Personal - Jaime y Natalia - 4192"
$ws.Range("A15").Value = "Cuenta Personal - Jaime y Nicolas - 4190"
$ws.Range("A16").Value = "Cuenta Personal - Natalia - 4229"
$ws.Range("A17").Value = "Cuenta Personal - Nicolas - 4230"
$ws.Range("A20").Value = "NNW Ventures LLC - 4357"
$ws.Range("A21").Value = "Cuenta Personal - Jacques - 4442"

# Update the sheet's active selection (view state) to A2:D21.
$ws.Activate()
$ws.Range("A2:D21").Select()
